# Logged Week 16 and performed season sim from Week 17
$wb = $excel.ActiveWorkbook

# --- OFF sheet: update "H" row target depth data ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 380
$wsOff.Range("C2").Value = 239
$wsOff.Range("D2").Value = 90
$wsOff.Range("E2").Value = 32
$wsOff.Range("F2").Value = 11

# --- DEF sheet: update "H" row target depth data ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 379
$wsDef.Range("C2").Value = 277
$wsDef.Range("D2").Value = 86
$wsDef.Range("E2").Value = 39
